# "goals working uwu + next scene"
#
# Tabelle1 ("timeLog"):
#   - E2:E19 turned into a single fill-down (shared) formula =D-C
#   - a new time-tracking entry was added in row 14 (Aris, 2024-02-01,
#     17:35 -> 19:00, "goalpoints + next level")
#   - the selection/scroll position moved to G13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Fill the duration formula down E2:E19 as one shared formula, exactly like
# selecting E2:E19 and typing Ctrl+Enter / dragging the fill handle from E2.
$ws.Range("E2:E19").Formula = "=D2-C2"

# Bring over row 13's cell formatting (so the new date cell B14 gets the
# same date number format as the rest of column B) before filling values.
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)

# New row of tracked time for Aris.
$ws.Range("A14").Value = "Aris"
$ws.Range("B14").Value = 45323
$ws.Range("C14").Formula = "=17+35/60"
$ws.Range("D14").Formula = "=19"
$ws.Range("F14").Value = "goalpoints + next level"

# Leave the selection where the author left it.
$ws.Range("G13").Select()
